# Update 'want-to-go count' (F) and 'minimum ticket price' (G) figures
# across all four worksheets to match the refreshed scrape data.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 328
$ws.Range("F3").Value = 531
$ws.Range("F4").Value = 621
$ws.Range("F6").Value = 456
$ws.Range("G6").Value = 78
$ws.Range("F9").Value = 856
$ws.Range("F10").Value = 820
$ws.Range("F11").Value = 400
$ws.Range("F12").Value = 71
$ws.Range("F13").Value = 421
$ws.Range("F16").Value = 884
$ws.Range("F18").Value = 29
$ws.Range("F19").Value = 1663
$ws.Range("F20").Value = 44
$ws.Range("F22").Value = 25
$ws.Range("F25").Value = 1453
$ws.Range("F27").Value = 519
$ws.Range("F28").Value = 348
$ws.Range("F29").Value = 576
$ws.Range("F30").Value = 408
$ws.Range("F31").Value = 2318
$ws.Range("F33").Value = 85
$ws.Range("F34").Value = 167
$ws.Range("F35").Value = 593
$ws.Range("F36").Value = 463
$ws.Range("F38").Value = 910
$ws.Range("F39").Value = 691
$ws.Range("F41").Value = 405
$ws.Range("F42").Value = 367

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 77
$ws.Range("F17").Value = 139
$ws.Range("G21").Value = 180
$ws.Range("F22").Value = 113
$ws.Range("F23").Value = 89
$ws.Range("F24").Value = 430

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 241
$ws.Range("F6").Value = 297

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 328
$ws.Range("F4").Value = 531
$ws.Range("F7").Value = 241
$ws.Range("F8").Value = 621
$ws.Range("F11").Value = 456
$ws.Range("G11").Value = 78
$ws.Range("F13").Value = 856
$ws.Range("F14").Value = 820
$ws.Range("F15").Value = 400
$ws.Range("F16").Value = 71
$ws.Range("F17").Value = 421
$ws.Range("F20").Value = 884
$ws.Range("F21").Value = 29
$ws.Range("F22").Value = 297
$ws.Range("F23").Value = 1663
$ws.Range("F24").Value = 44
$ws.Range("F31").Value = 1453
$ws.Range("F34").Value = 519
$ws.Range("F35").Value = 576
$ws.Range("F36").Value = 408
$ws.Range("F38").Value = 2318
$ws.Range("F39").Value = 85
$ws.Range("F40").Value = 167
$ws.Range("F41").Value = 593
$ws.Range("F42").Value = 463
$ws.Range("F44").Value = 910
$ws.Range("G46").Value = 180
$ws.Range("F47").Value = 430
$ws.Range("F48").Value = 691

